$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 3, pushing existing rows 3-22 down to 4-23
$ws.Rows.Item(3).Insert()

# Populate the newly inserted row 3 with the new weekly data point
$ws.Range("A3").Value = 11
$ws.Range("B3").Value = "Vega Monumental Concepción"
$ws.Range("C3").Value = "Bíobío"
$ws.Range("D3").Value = [DateTime]"2023-11-15"
$ws.Range("D3").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E3").Value = 8
$ws.Range("F3").Value = 300000000
$ws.Range("G3").Value = "Espárragos"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 300
$ws.Range("K3").Value = 1800
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = 1900
$ws.Range("N3").Value = "$/kilo"
$ws.Range("O3").Value = "Provincia de Linares"
$ws.Range("P3").Value = 1900
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = "Hortaliza"
